$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 5800
$ws.Cells.Item(40, 10).Value = 5800
$ws.Cells.Item(40, 12).Value = 5800
$ws.Cells.Item(40, 14).Value = -6150
$ws.Cells.Item(69, 8).Value = 21968.75
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 12).Value = 0
$ws.Cells.Item(69, 14).ClearContents()
$ws.Cells.Item(72, 8).Value = 21968.75
$ws.Cells.Item(72, 10).Value = 0
$ws.Cells.Item(72, 12).Value = 0
$ws.Cells.Item(72, 14).ClearContents()
$ws.Cells.Item(112, 8).Value = 2728.7827
$ws.Cells.Item(112, 9).Value = 1090.7142
$ws.Cells.Item(112, 11).Value = 3272.1426
$ws.Cells.Item(112, 13).Value = -2164.1426
$ws.Cells.Item(138, 8).Value = 3632.5312
$ws.Cells.Item(138, 9).Value = 1033.1818
$ws.Cells.Item(138, 11).Value = 3099.5454
$ws.Cells.Item(138, 13).Value = 2040.4546

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(24, 8).Value = 21338.5
$ws.Cells.Item(24, 10).Value = 21338.5
$ws.Cells.Item(24, 12).Value = 21338.5
$ws.Cells.Item(24, 14).Value = -22086.5
$ws.Cells.Item(45, 8).Value = 2883.3635
$ws.Cells.Item(45, 9).Value = 2318.2856
$ws.Cells.Item(45, 10).Value = 3872.25
$ws.Cells.Item(45, 11).Value = 2318.2856
$ws.Cells.Item(45, 12).Value = 3872.25
$ws.Cells.Item(45, 13).Value = -1941.2856
$ws.Cells.Item(45, 14).Value = -4626.25
$ws.Cells.Item(88, 8).Value = 1514.8572
$ws.Cells.Item(88, 9).Value = 1178.8889
$ws.Cells.Item(88, 10).Value = 2119.6
$ws.Cells.Item(88, 11).Value = 1178.8889
$ws.Cells.Item(88, 12).Value = 2119.6
$ws.Cells.Item(88, 13).Value = -772.8888999999999
$ws.Cells.Item(88, 14).Value = -2931.6
$ws.Cells.Item(91, 8).Value = 1514.8572
$ws.Cells.Item(91, 9).Value = 1178.8889
$ws.Cells.Item(91, 10).Value = 2119.6
$ws.Cells.Item(91, 11).Value = 1178.8889
$ws.Cells.Item(91, 12).Value = 2119.6
$ws.Cells.Item(91, 13).Value = 225.1111000000001
$ws.Cells.Item(91, 14).Value = -4927.6
$ws.Cells.Item(93, 8).Value = 30546
$ws.Cells.Item(93, 10).Value = 30546
$ws.Cells.Item(93, 12).Value = 30546
$ws.Cells.Item(93, 14).Value = -35538
$ws.Cells.Item(100, 8).Value = 21338.5
$ws.Cells.Item(100, 10).Value = 21338.5
$ws.Cells.Item(100, 12).Value = 21338.5
$ws.Cells.Item(100, 14).Value = -23502.5
$ws.Cells.Item(132, 8).Value = 1598.3611
$ws.Cells.Item(132, 9).Value = 1435.8889
$ws.Cells.Item(132, 11).Value = 4307.6667
$ws.Cells.Item(132, 13).Value = -1777.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 9005.833000000001
$ws.Cells.Item(20, 9).Value = 11683.875
$ws.Cells.Item(20, 10).Value = 3649.75
$ws.Cells.Item(20, 11).Value = 11683.875
$ws.Cells.Item(20, 12).Value = 3649.75
$ws.Cells.Item(20, 13).Value = -11436.875
$ws.Cells.Item(20, 14).Value = -4143.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2320.4
$ws.Cells.Item(31, 9).Value = 2183.1333
$ws.Cells.Item(31, 10).Value = 2732.2
$ws.Cells.Item(31, 11).Value = 2183.1333
$ws.Cells.Item(31, 12).Value = 2732.2
$ws.Cells.Item(31, 13).Value = -1888.1333
$ws.Cells.Item(31, 14).Value = -3322.2
$ws.Cells.Item(34, 8).Value = 2320.4
$ws.Cells.Item(34, 9).Value = 2183.1333
$ws.Cells.Item(34, 10).Value = 2732.2
$ws.Cells.Item(34, 11).Value = 2183.1333
$ws.Cells.Item(34, 12).Value = 2732.2
$ws.Cells.Item(34, 13).Value = -1981.1333
$ws.Cells.Item(34, 14).Value = -3136.2
$ws.Cells.Item(86, 8).Value = 2910.3333
$ws.Cells.Item(86, 9).Value = 2837
$ws.Cells.Item(86, 10).Value = 2947
$ws.Cells.Item(86, 11).Value = 2837
$ws.Cells.Item(86, 12).Value = 2947
$ws.Cells.Item(86, 13).Value = -1714
$ws.Cells.Item(86, 14).Value = -5193
$ws.Cells.Item(89, 8).Value = 2910.3333
$ws.Cells.Item(89, 9).Value = 2837
$ws.Cells.Item(89, 10).Value = 2947
$ws.Cells.Item(89, 11).Value = 14185
$ws.Cells.Item(89, 12).Value = 14735
$ws.Cells.Item(89, 13).Value = -8569
$ws.Cells.Item(89, 14).Value = -25967
$ws.Cells.Item(132, 8).Value = 2250
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 13).ClearContents()
$ws.Cells.Item(134, 8).Value = 6126
$ws.Cells.Item(134, 9).Value = 4875.0835
$ws.Cells.Item(134, 10).Value = 8627.833000000001
$ws.Cells.Item(134, 11).Value = 14625.2505
$ws.Cells.Item(134, 12).Value = 25883.499
$ws.Cells.Item(134, 13).Value = -12090.2505
$ws.Cells.Item(134, 14).Value = -30953.499
$ws.Cells.Item(141, 8).Value = 306247.62
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 10).Value = 306247.62
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 12).Value = 306247.62
$ws.Cells.Item(141, 13).ClearContents()
$ws.Cells.Item(141, 14).Value = -316607.62

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 38.5
$ws.Cells.Item(2, 10).Value = 32.375
$ws.Cells.Item(2, 12).Value = 194.25
$ws.Cells.Item(2, 14).Value = -420.25
$ws.Cells.Item(11, 8).Value = 3697.0417
$ws.Cells.Item(11, 9).Value = 4155.2383
$ws.Cells.Item(11, 11).Value = 12465.7149
$ws.Cells.Item(11, 13).Value = -12325.7149
$ws.Cells.Item(25, 8).Value = 3100
$ws.Cells.Item(25, 9).Value = 500
$ws.Cells.Item(25, 10).Value = 7000
$ws.Cells.Item(25, 11).Value = 1500
$ws.Cells.Item(25, 12).Value = 21000
$ws.Cells.Item(25, 13).Value = -1331
$ws.Cells.Item(25, 14).Value = -21338
$ws.Cells.Item(30, 8).Value = 3100
$ws.Cells.Item(30, 9).Value = 500
$ws.Cells.Item(30, 10).Value = 7000
$ws.Cells.Item(30, 11).Value = 1500
$ws.Cells.Item(30, 12).Value = 21000
$ws.Cells.Item(30, 13).Value = -1398
$ws.Cells.Item(30, 14).Value = -21204
$ws.Cells.Item(38, 8).Value = 199.90909
$ws.Cells.Item(38, 10).Value = 681.6667
$ws.Cells.Item(38, 12).Value = 2045.0001
$ws.Cells.Item(38, 14).Value = -2739.0001
$ws.Cells.Item(69, 8).Value = 750.5
$ws.Cells.Item(69, 9).Value = 750.5
$ws.Cells.Item(69, 11).Value = 2251.5
$ws.Cells.Item(69, 13).Value = -1440.5
$ws.Cells.Item(72, 8).Value = 750.5
$ws.Cells.Item(72, 9).Value = 750.5
$ws.Cells.Item(72, 11).Value = 6754.5
$ws.Cells.Item(72, 13).Value = -2698.5
$ws.Cells.Item(81, 8).Value = 10000
$ws.Cells.Item(81, 10).Value = 11000
$ws.Cells.Item(81, 12).Value = 33000
$ws.Cells.Item(81, 14).Value = -35246
$ws.Cells.Item(84, 8).Value = 10000
$ws.Cells.Item(84, 10).Value = 11000
$ws.Cells.Item(84, 12).Value = 99000
$ws.Cells.Item(84, 14).Value = -110232
$ws.Cells.Item(112, 8).Value = 2000
$ws.Cells.Item(112, 9).Value = 0
$ws.Cells.Item(112, 10).Value = 2000
$ws.Cells.Item(112, 11).Value = 0
$ws.Cells.Item(112, 12).Value = 6000
$ws.Cells.Item(112, 13).ClearContents()
$ws.Cells.Item(112, 14).Value = -8216
$ws.Cells.Item(118, 8).Value = 2013
$ws.Cells.Item(118, 10).Value = 2000
$ws.Cells.Item(118, 12).Value = 6000
$ws.Cells.Item(118, 14).Value = -8486
$ws.Cells.Item(131, 8).Value = 1594.8392
$ws.Cells.Item(131, 10).Value = 1615.2885
$ws.Cells.Item(131, 12).Value = 4845.8655
$ws.Cells.Item(131, 14).Value = -14925.8655
$ws.Cells.Item(141, 8).Value = 1695.7
$ws.Cells.Item(141, 9).Value = 1695.7
$ws.Cells.Item(141, 11).Value = 5087.1
$ws.Cells.Item(141, 13).Value = 92.89999999999964

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 1405.1471
$ws.Cells.Item(132, 9).Value = 1071.6154
$ws.Cells.Item(132, 10).Value = 2489.125
$ws.Cells.Item(132, 11).Value = 3214.8462
$ws.Cells.Item(132, 12).Value = 7467.375
$ws.Cells.Item(132, 13).Value = -684.8462
$ws.Cells.Item(132, 14).Value = -12527.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1795.7142
$ws.Cells.Item(16, 9).Value = 1795.7142
$ws.Cells.Item(16, 11).Value = 1795.7142
$ws.Cells.Item(16, 13).Value = -1625.7142
$ws.Cells.Item(74, 8).Value = 21158.2
$ws.Cells.Item(74, 9).Value = 21448
$ws.Cells.Item(74, 10).Value = 19999
$ws.Cells.Item(74, 11).Value = 21448
$ws.Cells.Item(74, 12).Value = 19999
$ws.Cells.Item(74, 13).Value = -20450
$ws.Cells.Item(74, 14).Value = -21995
$ws.Cells.Item(77, 8).Value = 21158.2
$ws.Cells.Item(77, 9).Value = 21448
$ws.Cells.Item(77, 10).Value = 19999
$ws.Cells.Item(77, 11).Value = 64344
$ws.Cells.Item(77, 12).Value = 59997
$ws.Cells.Item(77, 13).Value = -59352
$ws.Cells.Item(77, 14).Value = -69981
$ws.Cells.Item(93, 8).Value = 21929.842
$ws.Cells.Item(93, 9).Value = 1038.3846
$ws.Cells.Item(93, 11).Value = 1038.3846
$ws.Cells.Item(93, 13).Value = 209.6153999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 5288.8887
$ws.Cells.Item(4, 9).Value = 20500
$ws.Cells.Item(4, 11).Value = 20500
$ws.Cells.Item(4, 13).Value = -20387
$ws.Cells.Item(8, 8).Value = 0
$ws.Cells.Item(8, 9).Value = 0
$ws.Cells.Item(8, 11).Value = 0
$ws.Cells.Item(8, 13).ClearContents()
$ws.Cells.Item(9, 8).Value = 6
$ws.Cells.Item(9, 9).Value = 6
$ws.Cells.Item(9, 11).Value = 6
$ws.Cells.Item(9, 13).Value = 134
$ws.Cells.Item(11, 8).Value = 28747.5
$ws.Cells.Item(11, 10).Value = 28747.5
$ws.Cells.Item(11, 12).Value = 28747.5
$ws.Cells.Item(11, 14).Value = -29031.5
$ws.Cells.Item(13, 8).Value = 2473.25
$ws.Cells.Item(13, 9).Value = 1297.6666
$ws.Cells.Item(13, 11).Value = 1297.6666
$ws.Cells.Item(13, 13).Value = -1157.6666
$ws.Cells.Item(30, 8).Value = 24981
$ws.Cells.Item(30, 10).Value = 24981
$ws.Cells.Item(30, 12).Value = 24981
$ws.Cells.Item(30, 14).Value = -25195
$ws.Cells.Item(38, 8).Value = 11560
$ws.Cells.Item(38, 10).Value = 16059
$ws.Cells.Item(38, 12).Value = 16059
$ws.Cells.Item(38, 14).Value = -17005
$ws.Cells.Item(104, 8).Value = 8052
$ws.Cells.Item(104, 10).Value = 8052
$ws.Cells.Item(104, 12).Value = 8052
$ws.Cells.Item(104, 14).Value = -15040
$ws.Cells.Item(132, 8).Value = 1997.0588
$ws.Cells.Item(132, 9).Value = 2059.6875
$ws.Cells.Item(132, 11).Value = 6179.0625
$ws.Cells.Item(132, 13).Value = -3649.0625
$ws.Cells.Item(136, 8).Value = 3234.6667
$ws.Cells.Item(136, 9).Value = 3234.6667
$ws.Cells.Item(136, 11).Value = 9704.000100000001
$ws.Cells.Item(136, 13).Value = -7154.000100000001
